$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.028.08"
$ws.Range("E2").Value = "  -2.02%  "
$ws.Range("D3").Value = "3.372.68"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "571.30"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.53"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +8.68%  "
$ws.Range("D8").Value = "3.372.02"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.476"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.00%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.58"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +5.29%  "
$ws.Range("E11").Value = "  +3.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.388"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.22%  "
$ws.Range("D13").Value = "3.945.91"
$ws.Range("E13").Value = "  +0.09%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.121"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.98%  "
$ws.Range("E15").Value = "  +1.86%  "
$ws.Range("D16").Value = "3.368.77"
$ws.Range("E16").Value = "  -0.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.08"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +3.00%  "
$ws.Range("D18").Value = "61.175.16"
$ws.Range("E18").Value = "  -1.84%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.00"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +7.35%  "
$ws.Range("B20").Value = "Uniswap"
$ws.Range("C20").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "9.44"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.16%  "
$ws.Range("B21").Value = "Polkadot"
$ws.Range("C21").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.78"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "373.95"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.567"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +2.74%  "
$ws.Range("D24").Value = "3.507.17"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.56"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("E27").Value = "  +11.84%  "
$ws.Range("E28").Value = "  +22.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.67"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +12.28%  "
$ws.Range("E30").Value = "  +0.41%  "
$ws.Range("E31").Value = "  +4.57%  "
$ws.Range("E32").Value = "  +2.12%  "
$ws.Range("E33").Value = "  +4.43%  "
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "3.403.81"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "23.31"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +3.21%  "
$ws.Range("E37").Value = "  +9.30%  "
$ws.Range("B38").Value = "Aptos"
$ws.Range("C38").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.92"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +5.51%  "
$ws.Range("B39").Value = "ImmutableX"
$ws.Range("C39").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.56"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +6.11%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "163.26"
$ws.Range("D40").ClearFormats()
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0787"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.23%  "
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.41"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +4.59%  "
$ws.Range("E44").Value = "  +13.44%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.760"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("B46").Value = "OKB"
$ws.Range("C46").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "41.27"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("E47").Value = "  +4.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.11"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.74%  "
$ws.Range("E49").Value = "  +6.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "22.95"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +14.51%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.891"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +6.50%  "
